$d = $word.ActiveDocument

$d.Content.Find.Execute("81-37=44", $true, $false, $false, $false, $false, $true, 1, $false, "13+6=19", 2)
$d.Content.Find.Execute("53-32=21", $true, $false, $false, $false, $false, $true, 1, $false, "80+17=97", 2)
$d.Content.Find.Execute("90-9=81", $true, $false, $false, $false, $false, $true, 1, $false, "19+10=29", 2)
$d.Content.Find.Execute("49+9=58", $true, $false, $false, $false, $false, $true, 1, $false, "89-82=7", 2)
$d.Content.Find.Execute("1+56=57", $true, $false, $false, $false, $false, $true, 1, $false, "10+31=41", 2)
$d.Content.Find.Execute("46+53=99", $true, $false, $false, $false, $false, $true, 1, $false, "17+11=28", 2)
$d.Content.Find.Execute("59+8=67", $true, $false, $false, $false, $false, $true, 1, $false, "89-68=21", 2)
$d.Content.Find.Execute("26+43=69", $true, $false, $false, $false, $false, $true, 1, $false, "50-19=31", 2)
$d.Content.Find.Execute("45-21=24", $true, $false, $false, $false, $false, $true, 1, $false, "41-4=37", 2)
$d.Content.Find.Execute("68-52=16", $true, $false, $false, $false, $false, $true, 1, $false, "71-1=70", 2)
$d.Content.Find.Execute("43+20=63", $true, $false, $false, $false, $false, $true, 1, $false, "95-68=27", 2)
$d.Content.Find.Execute("27+47=74", $true, $false, $false, $false, $false, $true, 1, $false, "61+9=70", 2)
$d.Content.Find.Execute("9+82=91", $true, $false, $false, $false, $false, $true, 1, $false, "28-7=21", 2)
$d.Content.Find.Execute("71-8=63", $true, $false, $false, $false, $false, $true, 1, $false, "40+13=53", 2)
$d.Content.Find.Execute("91-73=18", $true, $false, $false, $false, $false, $true, 1, $false, "8-8=0", 2)
$d.Content.Find.Execute("53+29=82", $true, $false, $false, $false, $false, $true, 1, $false, "46-19=27", 2)
$d.Content.Find.Execute("52-50=2", $true, $false, $false, $false, $false, $true, 1, $false, "36-6=30", 2)
$d.Content.Find.Execute("26+19=45", $true, $false, $false, $false, $false, $true, 1, $false, "43+30=73", 2)
$d.Content.Find.Execute("31+3=34", $true, $false, $false, $false, $false, $true, 1, $false, "18+76=94", 2)
$d.Content.Find.Execute("93-7=86", $true, $false, $false, $false, $false, $true, 1, $false, "36+37=73", 2)
$d.Content.Find.Execute("39+52=91", $true, $false, $false, $false, $false, $true, 1, $false, "21+33=54", 2)
$d.Content.Find.Execute("70+6=76", $true, $false, $false, $false, $false, $true, 1, $false, "50-34=16", 2)
$d.Content.Find.Execute("21+17=38", $true, $false, $false, $false, $false, $true, 1, $false, "27+36=63", 2)
$d.Content.Find.Execute("46-21=25", $true, $false, $false, $false, $false, $true, 1, $false, "81-75=6", 2)
$d.Content.Find.Execute("52+40=92", $true, $false, $false, $false, $false, $true, 1, $false, "63+12=75", 2)
$d.Content.Find.Execute("81-49=32", $true, $false, $false, $false, $false, $true, 1, $false, "88-34=54", 2)
$d.Content.Find.Execute("55+23=78", $true, $false, $false, $false, $false, $true, 1, $false, "84-83=1", 2)
$d.Content.Find.Execute("66+33=99", $true, $false, $false, $false, $false, $true, 1, $false, "83-67=16", 2)
$d.Content.Find.Execute("94-81=13", $true, $false, $false, $false, $false, $true, 1, $false, "30+44=74", 2)
$d.Content.Find.Execute("84-9=75", $true, $false, $false, $false, $false, $true, 1, $false, "35-19=16", 2)
$d.Content.Find.Execute("18+15=33", $true, $false, $false, $false, $false, $true, 1, $false, "72-25=47", 2)
$d.Content.Find.Execute("13+10=23", $true, $false, $false, $false, $false, $true, 1, $false, "28+2=30", 2)
$d.Content.Find.Execute("78-70=8", $true, $false, $false, $false, $false, $true, 1, $false, "97-31=66", 2)
$d.Content.Find.Execute("49+23=72", $true, $false, $false, $false, $false, $true, 1, $false, "40-3=37", 2)
$d.Content.Find.Execute("0+54=54", $true, $false, $false, $false, $false, $true, 1, $false, "73-2=71", 2)
$d.Content.Find.Execute("37+31=68", $true, $false, $false, $false, $false, $true, 1, $false, "75-66=9", 2)
$d.Content.Find.Execute("36+40=76", $true, $false, $false, $false, $false, $true, 1, $false, "34+27=61", 2)
$d.Content.Find.Execute("87-83=4", $true, $false, $false, $false, $false, $true, 1, $false, "17+77=94", 2)
$d.Content.Find.Execute("15+83=98", $true, $false, $false, $false, $false, $true, 1, $false, "93-17=76", 2)
$d.Content.Find.Execute("61-57=4", $true, $false, $false, $false, $false, $true, 1, $false, "85-75=10", 2)
$d.Content.Find.Execute("59+1=60", $true, $false, $false, $false, $false, $true, 1, $false, "60-25=35", 2)
$d.Content.Find.Execute("27+24=51", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=67", 2)
$d.Content.Find.Execute("60-1=59", $true, $false, $false, $false, $false, $true, 1, $false, "90-57=33", 2)
$d.Content.Find.Execute("65-8=57", $true, $false, $false, $false, $false, $true, 1, $false, "33+22=55", 2)
$d.Content.Find.Execute("51-8=43", $true, $false, $false, $false, $false, $true, 1, $false, "24-4=20", 2)
$d.Content.Find.Execute("25+39=64", $true, $false, $false, $false, $false, $true, 1, $false, "43-40=3", 2)
$d.Content.Find.Execute("45-39=6", $true, $false, $false, $false, $false, $true, 1, $false, "61-7=54", 2)
$d.Content.Find.Execute("86-38=48", $true, $false, $false, $false, $false, $true, 1, $false, "53+21=74", 2)
$d.Content.Find.Execute("62-19=43", $true, $false, $false, $false, $false, $true, 1, $false, "93-88=5", 2)
$d.Content.Find.Execute("35+61=96", $true, $false, $false, $false, $false, $true, 1, $false, "95-86=9", 2)
$d.Content.Find.Execute("87-4=83", $true, $false, $false, $false, $false, $true, 1, $false, "63-12=51", 2)
$d.Content.Find.Execute("98-6=92", $true, $false, $false, $false, $false, $true, 1, $false, "36+55=91", 2)
$d.Content.Find.Execute("78-34=44", $true, $false, $false, $false, $false, $true, 1, $false, "50-26=24", 2)
$d.Content.Find.Execute("43+12=55", $true, $false, $false, $false, $false, $true, 1, $false, "27+27=54", 2)
$d.Content.Find.Execute("63-31=32", $true, $false, $false, $false, $false, $true, 1, $false, "35+48=83", 2)
$d.Content.Find.Execute("87-2=85", $true, $false, $false, $false, $false, $true, 1, $false, "80-46=34", 2)
$d.Content.Find.Execute("37-2=35", $true, $false, $false, $false, $false, $true, 1, $false, "97-8=89", 2)
$d.Content.Find.Execute("20+17=37", $true, $false, $false, $false, $false, $true, 1, $false, "49-42=7", 2)
$d.Content.Find.Execute("77-64=13", $true, $false, $false, $false, $false, $true, 1, $false, "47+39=86", 2)
$d.Content.Find.Execute("99-22=77", $true, $false, $false, $false, $false, $true, 1, $false, "30-29=1", 2)
$d.Content.Find.Execute("25+7=32", $true, $false, $false, $false, $false, $true, 1, $false, "27-21=6", 2)
$d.Content.Find.Execute("38+9=47", $true, $false, $false, $false, $false, $true, 1, $false, "66-27=39", 2)
$d.Content.Find.Execute("81-17=64", $true, $false, $false, $false, $false, $true, 1, $false, "88+11=99", 2)
$d.Content.Find.Execute("0+51=51", $true, $false, $false, $false, $false, $true, 1, $false, "35+9=44", 2)
$d.Content.Find.Execute("32+44=76", $true, $false, $false, $false, $false, $true, 1, $false, "1+86=87", 2)
$d.Content.Find.Execute("97-45=52", $true, $false, $false, $false, $false, $true, 1, $false, "11+72=83", 2)
$d.Content.Find.Execute("91-51=40", $true, $false, $false, $false, $false, $true, 1, $false, "86-22=64", 2)
$d.Content.Find.Execute("32-27=5", $true, $false, $false, $false, $false, $true, 1, $false, "1+74=75", 2)
$d.Content.Find.Execute("82-66=16", $true, $false, $false, $false, $false, $true, 1, $false, "81-39=42", 2)
$d.Content.Find.Execute("13+38=51", $true, $false, $false, $false, $false, $true, 1, $false, "12+42=54", 2)
$d.Content.Find.Execute("5+48=53", $true, $false, $false, $false, $false, $true, 1, $false, "12+30=42", 2)
$d.Content.Find.Execute("43+42=85", $true, $false, $false, $false, $false, $true, 1, $false, "88-29=59", 2)
$d.Content.Find.Execute("16+37=53", $true, $false, $false, $false, $false, $true, 1, $false, "21+41=62", 2)
$d.Content.Find.Execute("0+70=70", $true, $false, $false, $false, $false, $true, 1, $false, "39+53=92", 2)
$d.Content.Find.Execute("5+53=58", $true, $false, $false, $false, $false, $true, 1, $false, "55-34=21", 2)
$d.Content.Find.Execute("84-76=8", $true, $false, $false, $false, $false, $true, 1, $false, "53-26=27", 2)
$d.Content.Find.Execute("99-60=39", $true, $false, $false, $false, $false, $true, 1, $false, "96-18=78", 2)
$d.Content.Find.Execute("46-8=38", $true, $false, $false, $false, $false, $true, 1, $false, "50-41=9", 2)
$d.Content.Find.Execute("1+49=50", $true, $false, $false, $false, $false, $true, 1, $false, "86-56=30", 2)
$d.Content.Find.Execute("57-29=28", $true, $false, $false, $false, $false, $true, 1, $false, "23+20=43", 2)
$d.Content.Find.Execute("15+25=40", $true, $false, $false, $false, $false, $true, 1, $false, "32+18=50", 2)
$d.Content.Find.Execute("96-49=47", $true, $false, $false, $false, $false, $true, 1, $false, "15+67=82", 2)
$d.Content.Find.Execute("87-45=42", $true, $false, $false, $false, $false, $true, 1, $false, "92-66=26", 2)
$d.Content.Find.Execute("26+61=87", $true, $false, $false, $false, $false, $true, 1, $false, "67-62=5", 2)
$d.Content.Find.Execute("59+30=89", $true, $false, $false, $false, $false, $true, 1, $false, "44+23=67", 2)
$d.Content.Find.Execute("32+50=82", $true, $false, $false, $false, $false, $true, 1, $false, "95-73=22", 2)
$d.Content.Find.Execute("13-2=11", $true, $false, $false, $false, $false, $true, 1, $false, "56-42=14", 2)
$d.Content.Find.Execute("51+1=52", $true, $false, $false, $false, $false, $true, 1, $false, "9+35=44", 2)
$d.Content.Find.Execute("91-26=65", $true, $false, $false, $false, $false, $true, 1, $false, "72-33=39", 2)
$d.Content.Find.Execute("84-39=45", $true, $false, $false, $false, $false, $true, 1, $false, "67-64=3", 2)
$d.Content.Find.Execute("75-27=48", $true, $false, $false, $false, $false, $true, 1, $false, "92-19=73", 2)
$d.Content.Find.Execute("82-23=59", $true, $false, $false, $false, $false, $true, 1, $false, "68-0=68", 2)
$d.Content.Find.Execute("99-96=3", $true, $false, $false, $false, $false, $true, 1, $false, "23+50=73", 2)
$d.Content.Find.Execute("29+5=34", $true, $false, $false, $false, $false, $true, 1, $false, "30-14=16", 2)
$d.Content.Find.Execute("20+40=60", $true, $false, $false, $false, $false, $true, 1, $false, "35-3=32", 2)
$d.Content.Find.Execute("48+34=82", $true, $false, $false, $false, $false, $true, 1, $false, "35+48=83", 2)
$d.Content.Find.Execute("53+4=57", $true, $false, $false, $false, $false, $true, 1, $false, "58+1=59", 2)
$d.Content.Find.Execute("45+41=86", $true, $false, $false, $false, $false, $true, 1, $false, "2+17=19", 2)
$d.Content.Find.Execute("89-51=38", $true, $false, $false, $false, $false, $true, 1, $false, "91-2=89", 2)
$d.Content.Find.Execute("52+30=82", $true, $false, $false, $false, $false, $true, 1, $false, "6+44=50", 2)
